# Swap the order of "System" and the email address in column G
# wherever the recorded-by value is exactly "System, dnasr281@gmail.com".
# Other combinations (e.g. "backup@backdoor.com, System",
# "admin@admin.com, System", or a lone "dnasr281@gmail.com") are left
# untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G
    $current = $cell.Value2
    if ($current -eq $oldValue) {
        $cell.Value = $newValue
    }
}
